# Apply "Sideways Role cards done" edit:
# - Update the @Image column values for King/Traitor/Loyalist rows
#   from ".ai" to ".png" file extensions.
# - Update the active cell selection to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "King.png"
$ws.Range("B3").Value = "Traitor.png"
$ws.Range("B4").Value = "Loyalist.png"

$ws.Range("B3").Select()
